# Add a new "Alex Jordan" employee record to the Employee sheet (row 6),
# mirroring the existing rows, and hyperlink the password cell the same
# way the existing D5 ("johnSmith_@2023!!!") cell is linked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

$ws.Range("A6").Value = "Alex"
$ws.Range("B6").Value = "Jordan"
$ws.Range("C6").Value = "alexJordan"
$ws.Range("D6").Value = "jordan@_2023!!!"

# Mailto-style hyperlink on the new password cell, same pattern as D5.
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:jordan@_2023!!!")

# Match the existing hyperlink cell's style (the Hyperlink cell style).
$ws.Range("D6").Style = $ws.Range("D5").Style

# Reflect the saved selection/active cell recorded for the sheet.
[void]$ws.Range("G7").Select()
